$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix punctuation in proveedor names (comma -> period) ---
$ws.Range('E104').Value = 'FERNANDEZ. MARIO HUGO'
$ws.Range('E110').Value = 'ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN'
$ws.Range('F72').Value = 'MERCANZINI. GASTON ARIEL'

# --- Re-format "Importe" column from es-AR (1.234,56) to plain (1234.56) ---
$importeRange = $ws.Range("H2:H140")
$importeRange.NumberFormat = "@"

$ws.Range('H2').Value = '59800.00'
$ws.Range('H3').Value = '29835.00'
$ws.Range('H4').Value = '218000.00'
$ws.Range('H5').Value = '116.00'
$ws.Range('H6').Value = '12800.00'
$ws.Range('H7').Value = '100.00'
$ws.Range('H8').Value = '462.00'
$ws.Range('H9').Value = '4010.01'
$ws.Range('H10').Value = '27240.00'
$ws.Range('H11').Value = '622959.90'
$ws.Range('H12').Value = '17920.34'
$ws.Range('H13').Value = '2475.00'
$ws.Range('H14').Value = '11800.00'
$ws.Range('H15').Value = '2724.60'
$ws.Range('H16').Value = '22432.71'
$ws.Range('H17').Value = '2307.00'
$ws.Range('H18').Value = '17706.98'
$ws.Range('H19').Value = '6070.00'
$ws.Range('H20').Value = '1500.00'
$ws.Range('H21').Value = '29872.00'
$ws.Range('H22').Value = '6500.00'
$ws.Range('H23').Value = '1011.72'
$ws.Range('H24').Value = '16594.05'
$ws.Range('H25').Value = '1339.07'
$ws.Range('H26').Value = '3715.53'
$ws.Range('H27').Value = '12799.00'
$ws.Range('H28').Value = '2875.00'
$ws.Range('H29').Value = '5903.02'
$ws.Range('H30').Value = '750.01'
$ws.Range('H31').Value = '3700.00'
$ws.Range('H32').Value = '23207.50'
$ws.Range('H33').Value = '1225.00'
$ws.Range('H34').Value = '15092.00'
$ws.Range('H35').Value = '437.92'
$ws.Range('H36').Value = '97.00'
$ws.Range('H37').Value = '126000.00'
$ws.Range('H38').Value = '600.00'
$ws.Range('H39').Value = '439.02'
$ws.Range('H40').Value = '15947.00'
$ws.Range('H41').Value = '2222.00'
$ws.Range('H42').Value = '5796.00'
$ws.Range('H43').Value = '13894.70'
$ws.Range('H44').Value = '500.00'
$ws.Range('H45').Value = '4420.00'
$ws.Range('H46').Value = '22000.00'
$ws.Range('H47').Value = '55350.00'
$ws.Range('H48').Value = '1187.00'
$ws.Range('H49').Value = '4060.00'
$ws.Range('H50').Value = '469.70'
$ws.Range('H51').Value = '391178.31'
$ws.Range('H52').Value = '53261.49'
$ws.Range('H53').Value = '20000.00'
$ws.Range('H54').Value = '11.37'
$ws.Range('H55').Value = '68.56'
$ws.Range('H56').Value = '2950.00'
$ws.Range('H57').Value = '800.00'
$ws.Range('H58').Value = '5670.00'
$ws.Range('H59').Value = '4018.40'
$ws.Range('H60').Value = '419.00'
$ws.Range('H61').Value = '8770.00'
$ws.Range('H62').Value = '361.53'
$ws.Range('H63').Value = '512.00'
$ws.Range('H64').Value = '1000.00'
$ws.Range('H65').Value = '4874.00'
$ws.Range('H66').Value = '16852.50'
$ws.Range('H67').Value = '3500.00'
$ws.Range('H68').Value = '12000.00'
$ws.Range('H69').Value = '3200.00'
$ws.Range('H70').Value = '4000.00'
$ws.Range('H71').Value = '5570.00'
$ws.Range('H72').Value = '9000.00'
$ws.Range('H73').Value = '212100.00'
$ws.Range('H74').Value = '13500.00'
$ws.Range('H75').Value = '6100.00'
$ws.Range('H76').Value = '4734.00'
$ws.Range('H77').Value = '1205.00'
$ws.Range('H78').Value = '937.20'
$ws.Range('H79').Value = '219350.00'
$ws.Range('H80').Value = '37500.00'
$ws.Range('H81').Value = '11268.59'
$ws.Range('H82').Value = '2300.00'
$ws.Range('H83').Value = '1400.00'
$ws.Range('H84').Value = '4000.00'
$ws.Range('H85').Value = '1500.00'
$ws.Range('H86').Value = '72587.90'
$ws.Range('H87').Value = '3146.00'
$ws.Range('H88').Value = '1800.00'
$ws.Range('H89').Value = '1657.50'
$ws.Range('H90').Value = '1656.00'
$ws.Range('H91').Value = '2000.00'
$ws.Range('H92').Value = '1600.00'
$ws.Range('H93').Value = '2000.00'
$ws.Range('H94').Value = '6000.00'
$ws.Range('H95').Value = '950.00'
$ws.Range('H96').Value = '1500.00'
$ws.Range('H97').Value = '620.00'
$ws.Range('H98').Value = '4665.29'
$ws.Range('H99').Value = '6000.00'
$ws.Range('H100').Value = '1200.00'
$ws.Range('H101').Value = '1600.00'
$ws.Range('H102').Value = '8352.00'
$ws.Range('H103').Value = '1300.00'
$ws.Range('H104').Value = '5850.00'
$ws.Range('H105').Value = '10800.00'
$ws.Range('H106').Value = '668.41'
$ws.Range('H107').Value = '4200.00'
$ws.Range('H108').Value = '350.00'
$ws.Range('H109').Value = '3012.00'
$ws.Range('H110').Value = '1340.00'
$ws.Range('H111').Value = '4880.00'
$ws.Range('H112').Value = '1200.00'
$ws.Range('H113').Value = '200.00'
$ws.Range('H114').Value = '2100.00'
$ws.Range('H115').Value = '9629.00'
$ws.Range('H116').Value = '52.00'
$ws.Range('H117').Value = '351.00'
$ws.Range('H118').Value = '6505.00'
$ws.Range('H119').Value = '3949.50'
$ws.Range('H120').Value = '540.00'
$ws.Range('H121').Value = '4679.48'
$ws.Range('H122').Value = '6000.00'
$ws.Range('H123').Value = '10000.00'
$ws.Range('H124').Value = '1504177.80'
$ws.Range('H125').Value = '25000.00'
$ws.Range('H126').Value = '138780.00'
$ws.Range('H127').Value = '164169.00'
$ws.Range('H128').Value = '130000.00'
$ws.Range('H129').Value = '81000.00'
$ws.Range('H130').Value = '106870.00'
$ws.Range('H131').Value = '168782.26'
$ws.Range('H132').Value = '60000.00'
$ws.Range('H133').Value = '17500.00'
$ws.Range('H134').Value = '1174465.23'
$ws.Range('H135').Value = '11980.00'
$ws.Range('H136').Value = '27830.00'
$ws.Range('H137').Value = '527753.08'
$ws.Range('H138').Value = '20000.00'
$ws.Range('H139').Value = '938415.43'
$ws.Range('H140').Value = '70500.00'

$importeRange.ClearFormats()

Write-Output "done"
